$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 233040
$ws.Range("C2").Value = 281730
$ws.Range("D2").Value = 289092

$ws.Range("B3").Value = 233040
$ws.Range("C3").Value = 281730
$ws.Range("D3").Value = 289092

$ws.Range("B4").Value = 233040
$ws.Range("C4").Value = 258400
$ws.Range("D4").Value = 289092

$ws.Range("B5").Value = 230704
$ws.Range("C5").Value = 255914
$ws.Range("D5").Value = 289092

$ws.Range("B6").Value = 230704
$ws.Range("C6").Value = 247802
$ws.Range("D6").Value = 289092

$ws.Range("B7").Value = 230704
$ws.Range("C7").Value = 251448
$ws.Range("D7").Value = 289092

$ws.Range("B8").Value = 230704
$ws.Range("C8").Value = 256736
$ws.Range("D8").Value = 289092

$ws.Range("B9").Value = 230704
$ws.Range("C9").Value = 256278
$ws.Range("D9").Value = 289092

$ws.Range("B10").Value = 230704
$ws.Range("C10").Value = 257038
$ws.Range("D10").Value = 289092

$ws.Range("B11").Value = 232984
$ws.Range("C11").Value = 265920
$ws.Range("D11").Value = 289092

$ws.Range("B12").Value = 232164
$ws.Range("C12").Value = 263268
$ws.Range("D12").Value = 289092

$ws.Range("B13").Value = 224416
$ws.Range("C13").Value = 274820
$ws.Range("D13").Value = 289092

$ws.Range("B14").Value = 224416
$ws.Range("C14").Value = 274820
$ws.Range("D14").Value = 289092

$ws.Range("B15").Value = 229092
$ws.Range("C15").Value = 274820
$ws.Range("D15").Value = 289092

$ws.Range("B16").Value = 231902
$ws.Range("C16").Value = 274820
$ws.Range("D16").Value = 289092

$ws.Range("B17").Value = 224416
$ws.Range("C17").Value = 274820
$ws.Range("D17").Value = 289092

$ws.Range("B18").Value = 224416
$ws.Range("C18").Value = 274820
$ws.Range("D18").Value = 289092

$ws.Range("B19").Value = 230704
$ws.Range("C19").Value = 274820
$ws.Range("D19").Value = 289092

$ws.Range("B20").Value = 229982
$ws.Range("C20").Value = 274820
$ws.Range("D20").Value = 289092

$ws.Range("B21").Value = 224416
$ws.Range("C21").Value = 274820
$ws.Range("D21").Value = 289092

$ws.Range("B22").Value = 229982
$ws.Range("C22").Value = 274820
$ws.Range("D22").Value = 289092
